{"js": "const body = context.document.body;\nbody.paragraphs.load(\"items/text\");\nawait context.sync();\n\nconst paragraphs = body.paragraphs.items;\n\n// --- 1) Collapse the three \"CORE COMPETENCIES\" detail paragraphs into one ---\nlet coreIdx = -1;\nfor (let i = 0; i < paragraphs.length; i++) {\n  if (paragraphs[i].text.indexOf(\"Product Marketing Core: Market Intelligence\") === 0) {\n    coreIdx = i;\n    break;\n  }\n}\n\nif (coreIdx !== -1) {\n  const corePara = paragraphs[coreIdx];\n  corePara.insertText(\n    \"Product Marketing Core \u2022 Research & Analytics \u2022 Communication & Technology\",\n    Word.InsertLocation.replace\n  );\n  // Remove the two paragraphs that used to hold the detailed bullet lists.\n  paragraphs[coreIdx + 1].delete();\n  paragraphs[coreIdx + 2].delete();\n}\n\n// --- 2) Append a new \"TECHNICAL SKILLS\" section at the end of the document ---\nconst lastPara = paragraphs[paragraphs.length - 1];\n\nconst headingPara = lastPara.insertParagraph(\"TECHNICAL SKILLS\", Word.InsertLocation.after);\nconst productPara = headingPara.insertParagraph(\n  \"PRODUCT MARKETING CORE Market Intelligence & Competitive Analysis; Product Positioning & Messaging Development; Go-to-Market Strategy & Product Launch Management; Customer Segmentation & Buyer Persona Development; Cross-functional Team Leadership & Collaboration; Sales Enablement & Training Material Development; Data-Driven Decision Making & Analytics Interpretation\",\n  Word.InsertLocation.after\n);\nconst researchPara = productPara.insertParagraph(\n  \"RESEARCH & ANALYTICS Survey Methodology & Customer Insights; Market Research Design & Implementation; Competitive Intelligence & SWOT Analysis; Customer Journey Mapping & Behavioral Analysis; Statistical Modeling & Trend Analysis; Performance Metrics & Dashboard Development; A/B Testing & Conversion Optimization\",\n  Word.InsertLocation.after\n);\nconst commPara = researchPara.insertParagraph(\n  \"COMMUNICATION & TECHNOLOGY Strategic Messaging & Narrative Development; Technical Concept Translation for Business Audiences; Stakeholder Communication & Presentation Skills; Data Visualization & Reporting (Tableau, PowerBI, d3.js); Marketing Technology Stack Integration; Content Strategy & Thought Leadership; Client Relationship Management & Business Development\",\n  Word.InsertLocation.after\n);\n\n// Apply the heading style last so it doesn't bleed into the paragraphs\n// created after it while they were being chained together.\nheadingPara.style = \"Heading 2\";\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# --- 1) Collapse the three \"CORE COMPETENCIES\" detail paragraphs into one ---\n$coreIndex = -1\n$i = 0\nforeach ($p in $d.Paragraphs) {\n    $i = $i + 1\n    if ($p.Range.Text.StartsWith(\"Product Marketing Core: Market Intelligence\")) {\n        $coreIndex = $i\n        break\n    }\n}\n\nif ($coreIndex -gt 0) {\n    $corePar = $d.Paragraphs.Item($coreIndex)\n    $corePar.Range.Text = \"Product Marketing Core \u2022 Research & Analytics \u2022 Communication & Technology\"\n    # The two detailed paragraphs that used to follow are now right after the\n    # rewritten summary paragraph; delete them both.\n    $d.Paragraphs.Item($coreIndex + 1).Range.Delete()\n    $d.Paragraphs.Item($coreIndex + 1).Range.Delete()\n}\n\n# --- 2) Append a new \"TECHNICAL SKILLS\" section at the end of the document ---\n$count = $d.Paragraphs.Count\n$lastPar = $d.Paragraphs.Item($count)\n$endRange = $lastPar.Range\n$endRange.Collapse(0)\n$endRange.InsertAfter(\"`rTECHNICAL SKILLS`rPRODUCT MARKETING CORE Market Intelligence & Competitive Analysis; Product Positioning & Messaging Development; Go-to-Market Strategy & Product Launch Management; Customer Segmentation & Buyer Persona Development; Cross-functional Team Leadership & Collaboration; Sales Enablement & Training Material Development; Data-Driven Decision Making & Analytics Interpretation`rRESEARCH & ANALYTICS Survey Methodology & Customer Insights; Market Research Design & Implementation; Competitive Intelligence & SWOT Analysis; Customer Journey Mapping & Behavioral Analysis; Statistical Modeling & Trend Analysis; Performance Metrics & Dashboard Development; A/B Testing & Conversion Optimization`rCOMMUNICATION & TECHNOLOGY Strategic Messaging & Narrative Development; Technical Concept Translation for Business Audiences; Stakeholder Communication & Presentation Skills; Data Visualization & Reporting (Tableau, PowerBI, d3.js); Marketing Technology Stack Integration; Content Strategy & Thought Leadership; Client Relationship Management & Business Development\")\n\n$headingPar = $d.Paragraphs.Item($count + 1)\n$headingPar.Style = \"Heading 2\"\n"}
